$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has columns: A=Weapons, B=Enemy, C=Hp, D=Location, E=Item, F=Pet
# New layout needs two extra columns inserted:
#   - a new "Weapon_damage" column before Enemy (becomes column B, Enemy shifts to C)
#   - a new "Damage_dealt" column before Hp (becomes column D, Hp shifts to E)
# Using native column Insert() so formatting/styles shift exactly like Excel does.

$ws.Columns("B").Insert()
$ws.Columns("D").Insert()

# Header row (set D1 before B1 so the new shared strings are appended in the
# same order as the authored workbook: "Damage_dealt" then "Weapon_damage")
$ws.Range("D1").Value = "Damage_dealt"
$ws.Range("B1").Value = "Weapon_damage"

# Weapon_damage values (column B)
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 15
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 8
$ws.Range("B7").Value = 25
$ws.Range("B8").Value = 23
$ws.Range("B9").Value = 18
$ws.Range("B10").Value = 28
$ws.Range("B11").Value = 30

# Damage_dealt values (column D)
$ws.Range("D2").Value = 5
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 15
$ws.Range("D5").Value = 20
$ws.Range("D6").Value = 25
$ws.Range("D7").Value = 18
$ws.Range("D8").Value = 13
$ws.Range("D9").Value = 8
$ws.Range("D10").Value = 23
$ws.Range("D11").Value = 9

# Column widths for the two new columns (matches the authored widths of
# 15.90625 / 14.26953125 as closely as this engine's width rounding allows)
$ws.Columns("B").ColumnWidth = 14.92
$ws.Columns("D").ColumnWidth = 13.42

# Restore the selected cell like in the edited workbook
$ws.Range("C17").Select()
